# Generate Report for Handoff
# Updates the localization-status workbook after a new handoff xliff
# generation pass: refreshes the "Latest HO Xliff Generate Date" /
# "Latest Handoff Datetime" timestamps and bumps the "Priority" for the
# files that were re-handed-off from "low" to "ht".

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")
$overview = $wb.Worksheets.Item("Overview")

# Rows 4-7 on both locale sheets correspond to the four files that just
# had a fresh handoff xliff generated for them.
$zhcn.Range("E4:E7").Value = "ht"
$zhcn.Range("H4:H7").Value = "2016-08-30 08:34:59"

$dede.Range("E4:E7").Value = "ht"
$dede.Range("H4:H7").Value = "2016-08-30 08:35:15"

# Overview sheet mirrors the de-de "Latest HO Xliff Generate Date" value.
$overview.Range("G4:G7").Value = "2016-08-30 08:35:15"
